$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D/E hold plain text values (prices/volume-% are formatted
# strings, not numbers) in the source data. Temporarily mark the range
# as Text so Excel does not auto-coerce numeric-looking values (like
# "0.9997") into real numbers, then clear the format again so no
# lasting number-format change is left on the cells.
$dataRange = $ws.Range("D2:E51")
$dataRange.NumberFormat = "@"

$ws.Range('D2').Value = '29.038.42'
$ws.Range('E2').Value = '  +0.00%  '

$ws.Range('D3').Value = '1.828.76'
$ws.Range('E3').Value = '  -0.06%  '

$ws.Range('D4').Value = '0.9997'
$ws.Range('E4').Value = '  +0.11%  '

$ws.Range('D5').Value = '244.59'
$ws.Range('E5').Value = '  +1.41%  '

$ws.Range('D6').Value = '0.6312'
$ws.Range('E6').Value = '  +0.65%  '

$ws.Range('D7').Value = '0.9999'
$ws.Range('E7').Value = '  -0.01%  '

$ws.Range('D8').Value = '0.07537'
$ws.Range('E8').Value = '  -1.24%  '

$ws.Range('D9').Value = '0.2941'
$ws.Range('E9').Value = '  +0.86%  '

$ws.Range('D10').Value = '23.12'
$ws.Range('E10').Value = '  +1.43%  '

$ws.Range('D11').Value = '0.07703'
$ws.Range('E11').Value = '  +0.74%  '

$ws.Range('D12').Value = '1.828.51'
$ws.Range('E12').Value = '  -0.02%  '

$ws.Range('D13').Value = '4.999'
$ws.Range('E13').Value = '  +0.72%  '

$ws.Range('D14').Value = '0.6695'
$ws.Range('E14').Value = '  +0.61%  '

$ws.Range('D15').Value = '83.14'
$ws.Range('E15').Value = '  +0.83%  '

$ws.Range('D16').Value = '0.000009566'
$ws.Range('E16').Value = '  +1.96%  '

$ws.Range('D17').Value = '6.078'
$ws.Range('E17').Value = '  +1.48%  '

$ws.Range('D18').Value = '29.060.40'
$ws.Range('E18').Value = '  +0.72%  '

$ws.Range('D19').Value = '12.57'
$ws.Range('E19').Value = '  +1.95%  '

$ws.Range('D20').Value = '226.77'
$ws.Range('E20').Value = '  +0.78%  '

$ws.Range('D21').Value = '0.9991'
$ws.Range('E21').Value = '  -0.08%  '

$ws.Range('D22').Value = '7.143'
$ws.Range('E22').Value = '  -0.93%  '

$ws.Range('E23').Value = '  +0.00%  '

$ws.Range('D24').Value = '160.12'
$ws.Range('E24').Value = '  +0.12%  '

$ws.Range('D25').Value = '0.1427'
$ws.Range('E25').Value = '  +4.91%  '

$ws.Range('D26').Value = '8.509'
$ws.Range('E26').Value = '  +1.00%  '

$ws.Range('E27').Value = '  +0.72%  '

$ws.Range('E28').Value = '  +0.92%  '

$ws.Range('D29').Value = '4.145'
$ws.Range('E29').Value = '  +2.31%  '

$ws.Range('D30').Value = '4.071'
$ws.Range('E30').Value = '  +0.94%  '

$ws.Range('D31').Value = '0.05492'
$ws.Range('E31').Value = '  +5.50%  '

$ws.Range('E32').Value = '  -0.30%  '

$ws.Range('D33').Value = '1.858'
$ws.Range('E33').Value = '  +0.59%  '

$ws.Range('D34').Value = '0.7446'
$ws.Range('E34').Value = '  +1.81%  '

$ws.Range('D35').Value = '1.139'
$ws.Range('E35').Value = '  -1.33%  '

$ws.Range('D36').Value = '2.656'
$ws.Range('E36').Value = '  +1.68%  '

$ws.Range('D37').Value = '1.246.10'
$ws.Range('E37').Value = '  -2.16%  '

$ws.Range('D38').Value = '2.758'
$ws.Range('E38').Value = '  -0.03%  '

$ws.Range('D39').Value = '0.01784'
$ws.Range('E39').Value = '  -0.16%  '

$ws.Range('D40').Value = '6.589'
$ws.Range('E40').Value = '  +1.11%  '

$ws.Range('D41').Value = '0.9030'
$ws.Range('E41').Value = '  +1.43%  '

$ws.Range('E42').Value = '  -0.06%  '

$ws.Range('D43').Value = '101.41'
$ws.Range('E43').Value = '  -0.04%  '

$ws.Range('D44').Value = '1.979.60'
$ws.Range('E44').Value = '  +0.29%  '

$ws.Range('D45').Value = '65.02'
$ws.Range('E45').Value = '  +2.04%  '

$ws.Range('D46').Value = '0.00000000121'
$ws.Range('E46').Value = '  +1.10%  '

$ws.Range('D47').Value = '0.5102'
$ws.Range('E47').Value = '  -0.07%  '

$ws.Range('D48').Value = '0.4065'
$ws.Range('E48').Value = '  +2.19%  '

$ws.Range('D49').Value = '8.978'
$ws.Range('E49').Value = '  +1.46%  '

$ws.Range('D50').Value = '1.657'
$ws.Range('E50').Value = '  +0.74%  '

$ws.Range('B51').Value = 'Cronos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D51').Value = '0.05787'
$ws.Range('E51').Value = '  +0.87%  '

$dataRange.ClearFormats()